$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data of the paired rows (B/C order was reversed within each year group)
# Row 3 <-> Row 4
$ws.Range("A3").Value = "2016年C"
$ws.Range("B3").Value = 98.7
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 24
$ws.Range("E3").Value = 9293859.4
$ws.Range("A4").Value = "2016年B"
$ws.Range("B4").Value = 98.8
$ws.Range("C4").Value = 0.4
$ws.Range("D4").Value = 13.8
$ws.Range("E4").Value = 5909169.7

# Row 7 <-> Row 8
$ws.Range("A7").Value = "2017年C"
$ws.Range("B7").Value = 98.59999999999999
$ws.Range("C7").Value = -0.7
$ws.Range("D7").Value = 26.6
$ws.Range("E7").Value = 11333763.6
$ws.Range("A8").Value = "2017年B"
$ws.Range("B8").Value = 98.90000000000001
$ws.Range("C8").Value = 0.8
$ws.Range("D8").Value = 13.8
$ws.Range("E8").Value = 7355137.8

# Row 11 <-> Row 12
$ws.Range("A11").Value = "2018年C"
$ws.Range("B11").Value = 96.90000000000001
$ws.Range("C11").Value = -0.5
$ws.Range("D11").Value = 35.2
$ws.Range("E11").Value = 12485397.3
$ws.Range("A12").Value = "2018年B"
$ws.Range("B12").Value = 98.2
$ws.Range("C12").Value = -0.5
$ws.Range("D12").Value = 23.2
$ws.Range("E12").Value = 8333467.4

# Row 15 <-> Row 16
$ws.Range("A15").Value = "2019年C"
$ws.Range("B15").Value = 99.90000000000001
$ws.Range("C15").Value = 2.3
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 14268580.2
$ws.Range("A16").Value = "2019年B"
$ws.Range("B16").Value = 99
$ws.Range("C16").Value = 1.7
$ws.Range("D16").Value = 8.4
$ws.Range("E16").Value = 7772104.7

# Remove the F (产销率) and G (销售量) columns entirely
$ws.Range("F1:G17").Delete()
